# The document's first-page header carries the BTEC logo (originally
# embedded/exported as "image1.jpg") and the primary + first-page footers
# both carry the Pearson logo (originally exported as "image2.png").
# The authorial edit swaps these default picture names: the BTEC logo
# becomes "image2.jpg" while both Pearson logo instances become
# "image1.png". Rename every affected InlineShape accordingly.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# First-page header picture: BTec_Logo-Orange (image1.jpg -> image2.jpg)
$header = $sec.Headers.Item(2)
if ($header.Exists -and $header.Range.InlineShapes.Count -ge 1) {
    $headerLogo = $header.Range.InlineShapes.Item(1)
    $headerLogo.Name = "image2.jpg"
}

# First-page footer picture: Pearson logo (image2.png -> image1.png)
$firstFooter = $sec.Footers.Item(2)
if ($firstFooter.Exists -and $firstFooter.Range.InlineShapes.Count -ge 1) {
    $firstFooterLogo = $firstFooter.Range.InlineShapes.Item(1)
    $firstFooterLogo.Name = "image1.png"
}

# Default (primary) footer picture: Pearson logo (image2.png -> image1.png)
$defaultFooter = $sec.Footers.Item(1)
if ($defaultFooter.Exists -and $defaultFooter.Range.InlineShapes.Count -ge 1) {
    $defaultFooterLogo = $defaultFooter.Range.InlineShapes.Item(1)
    $defaultFooterLogo.Name = "image1.png"
}
